# Update weekly ranking [2026-01-14]
# Adds a new worksheet "magapoke_2026-01-14" at the end of the workbook,
# containing the rank/title ranking table for that week.

$titles = @(
    'ブルーロック',
    'WIND BREAKER',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    'ガチアクタ',
    '東京卍リベンジャーズ',
    'ベイビーステップ',
    'ギルティサークル',
    '島耕作',
    'イレギュラーズ',
    '黄昏町プリズナーズ',
    '魔女と傭兵',
    'ハードワーカー中田',
    '愛妻の裏アカ',
    '十字架のろくにん',
    '黒猫と魔女の教室',
    '南海トラフ巨大地震',
    'デッドアカウント',
    '君が僕らを悪魔と呼んだ頃',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    'となりの黒川さん',
    'ひゃくえむ。',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜',
    '異世界ウォーキング',
    '幼馴染とはラブコメにならない',
    'アルキメデスの大戦',
    'ドラハチ',
    'FAIRY TAIL 100 YEARS QUEST',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    'せいぶつ部の田辺くん',
    '食糧人類-Starving Anonymous-',
    '蒼く染めろ',
    'ともだちづくり',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'アオバノバスケ',
    'おやすみ ふみさん',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    'さわらないで小手指くん',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    'K-9~警視庁公安部公安第9課異能対策係~',
    'グラぱらっ！',
    '春くらり',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    'ハナバス　苔石花江のバスケ論',
    '普通の本はありません！',
    'ジュミドロ',
    'GALAXIAS',
    'いじめるヤバイ奴',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    'なれの果ての僕ら',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    '屋根の下のアルテミス',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    '我間乱 ―修羅―',
    'デスティニーラバーズ',
    'DAYS外伝',
    'ストーカー行為がバレて人生終了男',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    '可愛いだけじゃない式守さん',
    '降り積もれ孤独な死よ',
    '阿武ノーマル',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    '日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人',
    '剣帝学院の魔眼賢者',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    '「無能はいらない」と言われたから絶縁してやった　～最強の四天王に育てられた俺は、冒険者となり無双する～',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    '君が監督！',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    '鳴るさんだぁ',
    'ヒロインは絶望しました。',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    '復讐の教科書',
    '母という呪縛 娘という牢獄',
    '死ぬほど君の処女が欲しい',
    'MYS',
    'はっちぽっちぱんち',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    'インフルエンサーにストーカーされています',
    '人間消失',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    'イジらないで、長瀞さん',
    '魁の花巫女',
    'それがメイドのカンナです',
    '中華一番！極',
    '東京ネオンスキャンダル',
    'インフェクション',
    'ぼくのアデリア',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    '彼女、お借りします',
    'ウイニング パス'
)

$wb = $excel.ActiveWorkbook

$sheetName = "magapoke_2026-01-14"

# Use the previous week's sheet as a formatting template for the header row.
$templateSheet = $wb.Worksheets.Item("magapoke_2026-01-07")

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = $sheetName

# Match the outline summary settings used by the other weekly sheets
# (summaryBelow="1" summaryRight="1").
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match page margins used by the other weekly sheets (0.75/0.75/1/1in, 0.5/0.5in header/footer).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Copy the header row (rank/title), including its bold/centered/bordered style.
$templateSheet.Range("A1:B1").Copy($ws.Range("A1:B1"))

# Fill in the ranking rows (rank 1..100 in column A, title in column B),
# starting at row 2.
$row = 2
for ($i = 0; $i -lt $titles.Count; $i++) {
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
    $row = $row + 1
}
